# Apply "keep highest val strategy" updates to column D (pred) values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    12  = "AAA"
    15  = "BB"
    18  = "BB"
    30  = "AA"
    39  = "BB"
    51  = "B"
    59  = "BB"
    66  = "BBB"
    71  = "BB"
    76  = "BB"
    87  = "A"
    110 = "A"
    112 = "A"
    116 = "A"
    118 = "A"
    119 = "A"
    122 = "BB"
    129 = "A"
    130 = "A"
    134 = "A"
    138 = "BB"
    153 = "A"
    155 = "BB"
    158 = "BBB"
    159 = "A"
    161 = "AA"
    163 = "AA"
    180 = "AA"
    181 = "AA"
    182 = "BB"
    184 = "B"
    185 = "BB"
    193 = "BB"
    201 = "BB"
    212 = "BB"
    216 = "BB"
    226 = "BBB"
    232 = "AA"
    233 = "BB"
    237 = "B"
    239 = "AAA"
    265 = "BB"
    270 = "BBB"
    283 = "BB"
    285 = "BB"
    292 = "A"
    297 = "BBB"
    301 = "AA"
    302 = "AAA"
    306 = "A"
    308 = "BBB"
    311 = "A"
    313 = "B"
    323 = "B"
    342 = "BB"
    348 = "BB"
    370 = "BB"
    375 = "A"
    387 = "BBB"
    398 = "A"
}

foreach ($row in $updates.Keys) {
    $ws.Range("D$row").Value = $updates[$row]
}
